$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (prices and 1h volume % changes)
$updates = @{
    "D2" = "328.37"
    "E2" = "-0.69%"
    "D3" = "43.94"
    "E3" = "5.65%"
    "D4" = "5.465"
    "E4" = "-4.02%"
    "D5" = "0.08083"
    "E5" = "-3.67%"
    "D6" = "8.691"
    "E6" = "-1.44%"
    "D7" = "4.305"
    "E7" = "-3.80%"
    "D8" = "1.901"
    "E8" = "-5.24%"
    "D9" = "2.706"
    "E9" = "-7.90%"
    "D10" = "0.9409"
    "E10" = "1.71%"
    "D11" = "0.1210"
    "E11" = "-5.54%"
    "D12" = "0.1892"
    "E12" = "-4.53%"
    "D13" = "0.09548"
    "E13" = "0.34%"
    "D14" = "0.04152"
    "E14" = "7.88%"
    "D15" = "0.1072"
    "E15" = "0.97%"
    "D16" = "0.001286"
    "E16" = "-1.10%"
    "D17" = "0.006060"
    "E17" = "-0.78%"
    "D18" = "3.576"
    "E18" = "4.47%"
    "D20" = "8.495"
    "E20" = "-5.22%"
    "D21" = "0.1352"
    "E21" = "-0.80%"
    "D22" = "0.2604"
    "E22" = "3.77%"
    "D23" = "0.04369"
    "E23" = "-0.72%"
    "D24" = "0.001238"
    "E24" = "-2.79%"
    "D25" = "0.004298"
    "E25" = "-2.48%"
    "D26" = "0.0001234"
    "E26" = "3.61%"
    "D27" = "0.0004015"
    "E27" = "0.60%"
    "D39" = "0.02639"
    "E39" = "-8.02%"
    "D40" = "0.05446"
    "E40" = "-1.32%"
    "D41" = "0.007745"
    "E41" = "-2.70%"
    "D42" = "0.009760"
    "E42" = "8.42%"
    "E43" = "-2.98%"
    "D44" = "0.002127"
    "E44" = "2.70%"
    "D45" = "0.009895"
    "E45" = "-15.21%"
    "D46" = "0.00007327"
    "E46" = "5.94%"
    "D47" = "0.00000000755"
    "E47" = "0.61%"
    "D48" = "0.003556"
    "E48" = "2.67%"
    "D49" = "0.002285"
    "E49" = "0.27%"
    "D50" = "0.00002114"
    "E50" = "0.61%"
    "D51" = "0.0002013"
    "E51" = "0.61%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
